$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename "Diemen Sniep" to "Sniep" in cell A28
$ws.Range("A28").Value = "Sniep"

# Reflect the user's final selection on the sheet (moved from A23 to A28)
$ws.Range("A28").Select()
